# Apply the crypto-tracker refresh: updated Price (column D) and
# Volume(1h) (column E) values for rows 2-51, matching the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Column D/E cells in this sheet are stored as plain text (prices use
# "."-grouped thousands like "59.677.89", percentages keep their
# original padding like "  -3.74%  "). Excel auto-coerces a bare
# numeric-looking string (e.g. "552.78" or "7.00") into a real number
# when assigned to a General-formatted cell, which would silently
# change both the cell type and the displayed value (e.g. "7.00" ->
# "7"). To keep such values as text, the cell is switched to the
# "@" (Text) number format immediately before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, [string]$Addr, [string]$NewValue)

    $cell = $Sheet.Range($Addr)
    if ($NewValue -match '^[0-9]+(\.[0-9]+)?$') {
        # Force text storage so it round-trips as a string, not a float.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $NewValue
}

Set-CellText $ws "D2" "59.677.89"
Set-CellText $ws "E2" "  -3.74%  "
Set-CellText $ws "D3" "3.261.43"
Set-CellText $ws "E3" "  -4.58%  "
Set-CellText $ws "E4" "  +0.04%  "
Set-CellText $ws "D5" "552.78"
Set-CellText $ws "E5" "  -4.47%  "
Set-CellText $ws "D6" "140.18"
Set-CellText $ws "E6" "  -8.44%  "
Set-CellText $ws "E7" "  +0.08%  "
Set-CellText $ws "D8" "3.263.77"
Set-CellText $ws "E8" "  -4.50%  "
Set-CellText $ws "D9" "0.463"
Set-CellText $ws "E9" "  -4.22%  "
Set-CellText $ws "E10" "  -3.15%  "
Set-CellText $ws "D11" "0.117"
Set-CellText $ws "E11" "  -5.74%  "
Set-CellText $ws "D12" "0.401"
Set-CellText $ws "E12" "  -4.16%  "
Set-CellText $ws "D13" "3.830.21"
Set-CellText $ws "E13" "  -4.24%  "
Set-CellText $ws "E14" "  -0.48%  "
Set-CellText $ws "D15" "26.54"
Set-CellText $ws "E15" "  -7.31%  "
Set-CellText $ws "D16" "3.267.87"
Set-CellText $ws "E16" "  -4.53%  "
Set-CellText $ws "E17" "  -5.90%  "
Set-CellText $ws "D18" "59.810.59"
Set-CellText $ws "E18" "  -3.58%  "
Set-CellText $ws "D19" "6.03"
Set-CellText $ws "E19" "  -7.70%  "
Set-CellText $ws "D20" "13.62"
Set-CellText $ws "E20" "  -6.23%  "
Set-CellText $ws "D21" "8.42"
Set-CellText $ws "E21" "  -6.08%  "
Set-CellText $ws "D22" "371.34"
Set-CellText $ws "E22" "  -2.88%  "
Set-CellText $ws "D23" "72.65"
Set-CellText $ws "E23" "  -3.44%  "
Set-CellText $ws "E24" "  -0.10%  "
Set-CellText $ws "D25" "0.527"
Set-CellText $ws "E25" "  -7.76%  "
Set-CellText $ws "D26" "3.408.22"
Set-CellText $ws "E26" "  -4.27%  "
Set-CellText $ws "E27" "  -10.40%  "
Set-CellText $ws "E28" "  -4.93%  "
Set-CellText $ws "D29" "0.995"
Set-CellText $ws "E29" "  -0.34%  "
Set-CellText $ws "D30" "7.00"
Set-CellText $ws "E30" "  -8.89%  "
Set-CellText $ws "E31" "  -0.03%  "
Set-CellText $ws "D32" "2.00"
Set-CellText $ws "E32" "  -5.67%  "
Set-CellText $ws "D33" "7.41"
Set-CellText $ws "E33" "  -6.32%  "
Set-CellText $ws "D34" "22.34"
Set-CellText $ws "E34" "  -3.85%  "
Set-CellText $ws "D35" "1.22"
Set-CellText $ws "E35" "  -8.43%  "
Set-CellText $ws "D36" "166.26"
Set-CellText $ws "E36" "  -1.45%  "
Set-CellText $ws "D37" "5.02"
Set-CellText $ws "E37" "  -8.51%  "
Set-CellText $ws "E38" "  -5.57%  "
Set-CellText $ws "D39" "6.56"
Set-CellText $ws "E39" "  -5.78%  "
Set-CellText $ws "D40" "3.299.34"
Set-CellText $ws "E40" "  -4.34%  "
Set-CellText $ws "D41" "25.77"
Set-CellText $ws "E41" "  -16.53%  "
Set-CellText $ws "D42" "0.0718"
Set-CellText $ws "E42" "  -8.75%  "
Set-CellText $ws "D43" "41.26"
Set-CellText $ws "E43" "  -3.27%  "
Set-CellText $ws "D44" "0.741"
Set-CellText $ws "E44" "  -5.05%  "
Set-CellText $ws "D45" "4.07"
Set-CellText $ws "E45" "  -7.82%  "
Set-CellText $ws "E46" "  -4.64%  "
Set-CellText $ws "D47" "1.55"
Set-CellText $ws "E47" "  -7.64%  "
Set-CellText $ws "E48" "  +0.07%  "
Set-CellText $ws "D49" "2.305.15"
Set-CellText $ws "E49" "  -9.78%  "
Set-CellText $ws "D50" "6.29"
Set-CellText $ws "E50" "  -8.19%  "
Set-CellText $ws "D51" "21.01"
Set-CellText $ws "E51" "  -6.91%  "
